# PROYECTO_CS_10_08_ CO.docx -- wording/spelling fixes ("primera edición")
#
# The underlying XML diff shows Word re-inserting <w:proofErr/> spell/grammar
# check markers (a side-effect of opening + re-saving in a newer Word that
# re-ran the proofer) around many runs; those are not semantically
# meaningful edits on their own. The genuine content changes are a handful
# of wording/spelling corrections scattered through the body text. We apply
# those with Find/Replace, which is the operation a human editor actually
# performed.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $find"
    }
}

# 1. "Proyecto" lower-cased, comma after "razón", "2 o 3" spelled out, comma after "forma"
Replace-Text "Este Proyecto requiere ser realizado en equipo. Por tal razón debes reunirte con 2 o 3 de tus compañeros de clase y distribuir las tareas que se asignarán, la idea es que sea un trabajo que permita la discusión y el debate. De igual forma es una buena oportunidad para construir colectivamente una propuesta para presentarle al profesor y al resto del grupo. " `
              "Este proyecto requiere ser realizado en equipo. Por tal razón, debes reunirte con dos o tres de tus compañeros de clase y distribuir las tareas que se asignarán, la idea es que sea un trabajo que permita la discusión y el debate. De igual forma, es una buena oportunidad para construir colectivamente una propuesta para presentarle al profesor y al resto del grupo. "

# 2. "configurar" -> "comprender"
Replace-Text "va a servir para configurar mejor los cambios" `
              "va a servir para comprender mejor los cambios"

# 3. "deben" -> "debes", add "con tus compañeros", accent on "cuáles", drop comma, "realizar" -> "hacerles", "5" -> "cinco"
Replace-Text "En grupo deben debatir cuales son las personas más antiguas de sus familias, a quienes puedan realizar una entrevista. La idea es tener 5 fuentes de información." `
              "En grupo debes debatir con tus compañeros cuáles son las personas más antiguas de sus familias a quienes puedan hacerles una entrevista. La idea es tener cinco fuentes de información."

# 4. "5" -> "cinco", "realiza" -> "redactar", fix double space, "creencias, razas" -> "creencias y razas", "poder realizar" -> "hacer"
Replace-Text "Una vez escogidas las 5 personas a entrevistar, deben realiza un cuestionario que permita describir  la procedencia de las familias, las costumbres, creencias, razas. El objetivo es poder realizar un ejercicio comparativo con lo que sucede hoy. " `
              "Una vez escogidas las cinco personas a entrevistar, deben redactar un cuestionario que permita describir la procedencia de las familias, las costumbres, creencias y razas. El objetivo es hacer un ejercicio comparativo con lo que sucede hoy. "

# 5. "5" -> "cinco", accent on "cuáles"
Replace-Text "Luego de tener las 5 entrevistas, deben escucharlas grupalmente y definir cuales son los cambios y continuidades de las prácticas escuchadas. " `
              "Luego de tener las cinco entrevistas, deben escucharlas grupalmente y definir cuáles son los cambios y continuidades de las prácticas escuchadas. "

# 6. add comma after "anterior", "exponer al" -> "exponerle al", drop comma before "los hallaz..."
Replace-Text "Terminado el debate del punto anterior deben proponerse exponer al" `
              "Terminado el debate del punto anterior, deben proponerse exponerle al"
Replace-Text " grupo y al profesor, los hallaz" `
              " grupo y al profesor los hallaz"

# 7. turn the list into a parenthetical and add "emplear"
Replace-Text "La presentación deben hacerla de la forma más creativa posible, mímica, un audio, un performance, la metodología que les permita contar de mejor forma lo analizado. " `
              "La presentación deben hacerla de la forma más creativa posible (mímica, un audio, un performance), emplear la metodología que les permita contar de mejor forma lo analizado.  "

# 8. fix accidental double space
Replace-Text "que propone un trabajo grupal  de investigación sobre las procedencias de las personas de su entorno." `
              "que propone un trabajo grupal de investigación sobre las procedencias de las personas de su entorno."

# 9. drop the now-redundant lone-space run that trailed the "_GoBack" bookmark
#    (its single space is already absorbed into the merged run from step 7,
#    which now ends "...analizado.  " with two trailing spaces).
$bm = $d.Bookmarks("_GoBack")
$bmRange = $bm.Range
$paraRange = $bmRange.Paragraphs(1).Range
$tail = $d.Range($bmRange.End, $paraRange.End)
if ($tail.Text -eq " `r") {
    # keep the paragraph mark, drop only the stray space before it
    $spaceOnly = $d.Range($bmRange.End, $paraRange.End - 1)
    $spaceOnly.Text = ""
}

Write-Output "done"
